$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data columns (rows 2-25):
# I: 1 -> 2, K: 2 -> 1, M: 1 -> 2, O: 2 -> 1
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# Add new header cells P1/Q1 (continuing the numeric sequence)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header formatting (bold, bordered, centered) from an existing
# header cell onto the new header cells.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Add new data columns P and Q (all rows 2-25 get value 2)
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2

$excel.CutCopyMode = 0
